# Example.xlsx: add a new "Strength Value" column (E) mirroring the
# existing "Strength" column (C) — fixture data for the fieldName /
# PropertyName mapping change described in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Strength Value"
$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 5
$ws.Range("E4").Value = 7
$ws.Range("E5").Value = 9
$ws.Range("E6").Value = 11
$ws.Range("E7").Value = 13
$ws.Range("E8").Value = 15
$ws.Range("E9").Value = 17
$ws.Range("E10").Value = 19
$ws.Range("E11").Value = 21

# Cosmetic follow-up from the same save (widen the new column, move the
# active selection, zoom in, and set the page setup as in the workbook).
$ws.Columns.Item(5).ColumnWidth = 12.14

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("G13").Select()
$excel.ActiveWindow.Zoom = 200
